# Add extra columns (WIN, TOP4, TOP5, TOP6, RELEGATION) before ExpPoints,
# and refresh team order / ExpPoints values for the matchday-2 prediction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: Rank | Team | WIN | TOP4 | TOP5 | TOP6 | RELEGATION | ExpPoints
# Move the old "ExpPoints" header from C1 to H1, and insert the new headers
# in between, keeping the existing header formatting (bold, thin border,
# centered / top-aligned) used by A1:B1 - copy the format from an existing
# header cell rather than re-building it, so the same style is reused.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

$ws.Range("B1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Updated table (rank, team, expected points) in final row order.
$teams = @(
    @(1,  "Barcelona",            84.14071792858097),
    @(2,  "Real Madrid",          83.38857605388084),
    @(3,  "Atlético de Madrid",   70.84319889899569),
    @(4,  "Villarreal",           64.46129126634317),
    @(5,  "Real Betis",           58.89738184341011),
    @(6,  "Athletic Club",        55.80746759644336),
    @(7,  "Rayo Vallecano",       49.92755595820069),
    @(8,  "Celta de Vigo",        49.0576433883055),
    @(9,  "Osasuna",              47.8695611366082),
    @(10, "Getafe",               47.59442515289825),
    @(11, "Espanyol",             46.79836602063359),
    @(12, "Real Sociedad",        45.58760964599742),
    @(13, "Alavés",               44.42944581169475),
    @(14, "Valencia",             43.26029905148187),
    @(15, "Sevilla",              42.33006931391447),
    @(16, "Mallorca",             42.28582347080862),
    @(17, "Elche",                42.08959791719842),
    @(18, "Levante",              37.67572957516552),
    @(19, "Girona",               34.7148744783816),
    @(20, "Real Oviedo",          33.70365341973389)
)

for ($i = 0; $i -lt $teams.Count; $i++) {
    $row = $i + 2
    $rank = $teams[$i][0]
    $team = $teams[$i][1]
    $pts  = $teams[$i][2]

    $ws.Cells.Item($row, 1).Value = $rank
    $ws.Cells.Item($row, 2).Value = $team
    # Clear the new intermediate columns - reserved for future Monte Carlo
    # simulation percentages (WIN / TOP4 / TOP5 / TOP6 / RELEGATION).
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = $pts
}

Write-Output "done"
